$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('zoxmqkn@xvimixw.com', 'bRf3dxR', 'rmwidji', 'Travell')
    ,@('opfkeurslv@oxomotqawj.com', 'VqN2dT6WNP', 'fpvjiyvtvt', 'Lifelong l')
    ,@('ktqtlhrpdyn@xjjeqjouifr.com', 'YmfV8xrXKbA', 'xiuzmpaiojn', 'I love to c')
    ,@('fcoxxtcldvb@owrppwzrswo.com', 'bRIc8Ep2Chr', 'rlovzeosbsk', 'Software en')
    ,@('pejsbchtmte@rthrommxpib.com', 'uR6uEok316Y', 'dqeguqmndih', 'Traveller a')
    ,@('ymiumofqgsq@eiesgtbysxi.com', '4u0NQcTs5JJ', 'kadoejagytz', 'I love to c')
    ,@('kxkdhxomuyx@jpvkcsmxwkp.com', 'JuWLBXrBq59', 'fxvamidadsl', 'I love to c')
    ,@('ltwkmwglwi@semnuldakb.com', 'H1iG8OylMx', 'rzspuoiwfj', 'Traveller ')
    ,@('tbskknfhe@rddfpzsyi.com', 'xiBAi6iz9', 'gvbmwvjjj', 'Software ')
    ,@('pyohgqw@vsteqjb.com', '8FdPJEn', 'inwewdn', 'Foodie.')
    ,@('tczjmm@vqsmlp.com', 'KZ2hrP', 'qzogiv', 'Travel')
    ,@('ziojaz@mvtoif.com', '3IREI7', 'ibuitw', 'Foodie')
    ,@('cyppyg@reeurc.com', 'NCBKV6', 'coipxl', 'Travel')
    ,@('ctrwlna@ztvzuum.com', 'Hpl6lfl', 'rvgyvnj', 'Foodie.')
    ,@('bncdykcxev@cbqchvwsfl.com', '4cgoBdtIdd', 'baalnaprgn', 'Lifelong l')
    ,@('fcazdejoyy@zjetbmkqpg.com', 'BnAQRfGnJT', 'obgrgsfgbo', 'Traveller ')
    ,@('pzfsoa@ctqcyo.com', '1kDEdI', 'fwgmed', 'Softwa')
    ,@('iwksmntrjn@ffyimoytpk.com', 'yu2uztCXAm', 'shwbnpkses', 'Foodie. Yo')
    ,@('sjvaunm@muhlsoz.com', 'c1AzmWM', 'cnqzxkb', 'I love ')
    ,@('tautpgrr@bulnrzkn.com', '1ifoGNNG', 'jzgnutfs', 'Foodie. ')
    ,@('jcgzholj@usipwimf.com', 'Qs8AeXjY', 'rfavrkxv', 'Travelle')
    ,@('srmontf@sdlepvq.com', 'hpy7S4n', 'rtaxaxs', 'Travell')
    ,@('cjgjwmepkj@qokiqlxynf.com', 'ep0W3pdb86', 'mdjmmeqmiw', 'Foodie. Yo')
    ,@('vncwqxmx@ejkwouwr.com', '8BLzUZn3', 'csirimsz', 'Travelle')
    ,@('gxnbqevc@flyrxuwa.com', 'RPGQnPQT', 'xayqiosg', 'Travelle')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
